# Edit function to clear terminal window
# Add a new user row (Irina / irina@gmail.com) to the Users sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "Irina"
$ws.Range("B9").Value = "irina@gmail.com"
